$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 22: Do Xuan Thanh ----
$ws.Range("A22").Value = "9edbd47a-186b-4c5b-a077-8446b7418f6f"
$ws.Range("B22").Value = "xuanthanh"
$ws.Range("C22").Value = "Customer_Service"
$ws.Range("D22").Value = "xuanthanh@gmail.com"
$ws.Range("E22").Value = "'0346494851"
$ws.Range("F22").Value = '$2y$10$Sgk31oRzu3kNLfIUvAJLNOaTwts89qxic3Yzk2s59o0C/yMxabh2G'
$ws.Range("G22").Value = "Đỗ Xuân Thanh"
$ws.Range("H22").Value = "https://vapa.vn/wp-content/uploads/2022/12/anh-avatar-facebook-dep-001.jpg"
$ws.Range("I22").Value = "Male"
# J22 (dob) must stay a plain text value identical to the existing J20 cell;
# copy/paste-values from that cell avoids Excel's automatic text->date coercion.
$ws.Range("J20").Copy()
$ws.Range("J22").PasteSpecial(-4163)
$ws.Range("K22").Value = "Bình Dương"
$ws.Range("L22").Value = $true
$ws.Range("M22").Value = $true
$ws.Range("N22").Value = "Đại học Y khoa Vinh"

# ---- Row 23: Le Ngoc Nhu ----
$ws.Range("A23").Value = "2798c948-07a7-4f85-b7a6-8d8d69e53676"
$ws.Range("B23").Value = "ngocnhu"
$ws.Range("C23").Value = "Customer_Service"
$ws.Range("D23").Value = "ngocnhu@gmail.com"
$ws.Range("E23").Value = "'0366995813"
$ws.Range("F23").Value = '$2y$10$Sgk31oRzu3kNLfIUvAJLNOaTwts89qxic3Yzk2s59o0C/yMxabh2G'
$ws.Range("G23").Value = "Lê Ngọc Như"
$ws.Range("H23").Value = "https://thao68.com/wp-content/uploads/2022/03/avatar-facebook-3.jpg"
$ws.Range("I23").Value = "Female"
# J23 (dob) copied from the existing J10 cell for the same reason as above.
$ws.Range("J10").Copy()
$ws.Range("J23").PasteSpecial(-4163)
$ws.Range("K23").Value = "Q8, TP.HCM"
$ws.Range("L23").Value = $true
$ws.Range("M23").Value = $true
$ws.Range("N23").Value = "Đại học Y khoa Vinh"

$ws.Range("C23").Select()
